# CATATAN REVISI.xlsx - "cetak keseluruhan penambahan CATATAN REVISI.xlsx perbaikan tampilan"
#
# Adds 6 new chat-log rows (92-97) to the bottom of the CATATAN/STATUS table on
# Sheet1, extends the table/autofilter range to cover them, and moves the
# active selection down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Copy the existing row formatting down onto the new rows -----------
# Rows 92-95 use the same look as the earlier "CATATAN" rows (style used on
# row 3: wrapped text, borderId with bottom edge).
$ws.Range("B3:C3").Copy() | Out-Null
$ws.Range("B92:C95").PasteSpecial(-4122) | Out-Null

# Rows 96-97 use the look of the trailing rows (style used on row 91:
# wrapped text, borderId without bottom edge, since they sit at the bottom
# of the table).
$ws.Range("B91:C91").Copy() | Out-Null
$ws.Range("B96:C97").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 2. Fill in the new CATATAN text (STATUS column stays blank) ----------
$ws.Cells.Item(92, 2).Value() = "[5/13/2015, 00:41] Pak Devi Leuwigajah: Za barusan udah dikirim data yang hasil pembangunan sama yang akan di bangun 2015. Cobain aja dulu tar klo ada yanh aneh peta nya kasi tau"
$ws.Cells.Item(93, 2).Value() = "[5/13/2015, 00:42] Pak Devi Leuwigajah: Itu kelemahan peta nya klo ada belokan susah jadi cuma bisa lurus aja"
$ws.Cells.Item(94, 2).Value() = "[5/13/2015, 00:42] Pak Devi Leuwigajah: Untuk jalan sama drainase"
$ws.Cells.Item(95, 2).Value() = "[5/13/2015, 18:09] Pak Devi Leuwigajah: Za gimana kemaren udah bisa di masukin data base yang di email nya?"
$ws.Cells.Item(96, 2).Value() = "[5/13/2015, 18:09] Pak Devi Leuwigajah: Ada yang aneh ga koordinat nya"
$ws.Cells.Item(97, 2).Value() = "Terus za ini data usulan dah mulai ada. Yang di masukinnya nanti di exel berdasrkan kegiatan aja jadi dari rw.01 sampe 20"

# Row 92 holds a long note and is taller than a single text line, like the
# other two-line rows already in the sheet (e.g. row 87).
$ws.Rows.Item(92).RowHeight = 30

# --- 3. Grow the Table1 / autofilter range so it covers the new rows ------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B1:C97"))

# --- 4. Move the selection to the new last row, matching the edited file --
$ws.Range("B97").Select() | Out-Null
